$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a brand new "2022-Q3" worksheet right after "总计" and
#    before the existing "2022-Q2" sheet.
# ------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)

$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# Re-fetch the "2022-Q2" sheet *after* inserting the new sheet: sheet
# references resolve by live position, and inserting a sheet shifts
# every later index, so a reference grabbed beforehand would now point
# at the wrong tab.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Clone the header/row-index formatting from the "2022-Q2" sheet so the
# new sheet's styling (bold header row, bordered index column, etc.)
# matches the rest of the workbook instead of Excel's defaults.
$q2Sheet.Range("A1:H3").Copy()
$q3Sheet.Range("A1:H3").PasteSpecial(-4122)

# Headers
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Row 2: 浙商智选食品饮料股票A
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'011179"
$q3Sheet.Range("C2").Value = "浙商智选食品饮料股票A"
$q3Sheet.Range("D2").Value = "'0.14"
$q3Sheet.Range("E2").Value = "'91.42"
$q3Sheet.Range("F2").Value = "'5.49"
$q3Sheet.Range("G2").Value = "'0.0077"
$q3Sheet.Range("H2").Value = 10

# Row 3: 浙商智选食品饮料股票C
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'011180"
$q3Sheet.Range("C3").Value = "浙商智选食品饮料股票C"
$q3Sheet.Range("D3").Value = "'0.08"
$q3Sheet.Range("E3").Value = "'91.42"
$q3Sheet.Range("F3").Value = "'5.49"
$q3Sheet.Range("G3").Value = "'0.0044"
$q3Sheet.Range("H3").Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: the old "2022-Q2" row becomes
#    "2022-Q3", a new row is inserted for "2022-Q2" (re-using the old
#    counts/values), and the "2021-Q4" row shifts down one row with
#    its running index bumped.
# ------------------------------------------------------------------
$summarySheet.Range("B2").Value = "2022-Q3"

$summarySheet.Rows.Item(3).Insert()

# Re-create the bold/bordered index-column style on the newly inserted
# row's A cell (matches the style already used by A2/A4).
$summarySheet.Range("A3").Font.Bold = $true
$summarySheet.Range("A3").HorizontalAlignment = -4108
$summarySheet.Range("A3").VerticalAlignment = -4160
$summarySheet.Range("A3").Borders.LineStyle = 1

$summarySheet.Range("A3").Value = 1
$summarySheet.Range("B3").Value = "2022-Q2"
$summarySheet.Range("C3").Value = 2
$summarySheet.Range("D3").Value = 0.01

$summarySheet.Range("A4").Value = 2

Write-Host "2022-Q3 sheet added and summary sheet updated"
